$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 5856
$ws.Range("B2").Value = 5461
$ws.Range("C2").Value = 7904
$ws.Range("D2").Value = 2338
$ws.Range("E2").Value = 5988
$ws.Range("F2").Value = 3757
